$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at position 328 (pushes existing rows 328.. down to 330..)
$ws.Rows("328:329").Insert()

# Populate the two newly-inserted rows with the new weekly price data.
# Row 328: "Primera" quality
$ws.Cells.Item(328, 1).Value2  = 11
$ws.Cells.Item(328, 2).Value   = "Vega Monumental Concepción"
$ws.Cells.Item(328, 3).Value   = "Bíobío"
$ws.Cells.Item(328, 4).Value2  = 44992
$ws.Cells.Item(328, 5).Value2  = 8
$ws.Cells.Item(328, 6).Value2  = 100114014
$ws.Cells.Item(328, 7).Value   = "Betarraga"
$ws.Cells.Item(328, 8).Value   = "Sin especificar"
$ws.Cells.Item(328, 9).Value   = "Primera"
$ws.Cells.Item(328, 10).Value2 = 500
$ws.Cells.Item(328, 11).Value2 = 600
$ws.Cells.Item(328, 12).Value2 = 700
$ws.Cells.Item(328, 13).Value2 = 660
$ws.Cells.Item(328, 14).Value  = "$/paquete 5 unidades"
$ws.Cells.Item(328, 15).Value  = "Región Metropolitana"
$ws.Cells.Item(328, 16).Value2 = 132
$ws.Cells.Item(328, 17).Value2 = 5
$ws.Cells.Item(328, 18).Value  = "Hortaliza"

# Row 329: "Segunda" quality
$ws.Cells.Item(329, 1).Value2  = 11
$ws.Cells.Item(329, 2).Value   = "Vega Monumental Concepción"
$ws.Cells.Item(329, 3).Value   = "Bíobío"
$ws.Cells.Item(329, 4).Value2  = 44992
$ws.Cells.Item(329, 5).Value2  = 8
$ws.Cells.Item(329, 6).Value2  = 100114014
$ws.Cells.Item(329, 7).Value   = "Betarraga"
$ws.Cells.Item(329, 8).Value   = "Sin especificar"
$ws.Cells.Item(329, 9).Value   = "Segunda"
$ws.Cells.Item(329, 10).Value2 = 300
$ws.Cells.Item(329, 11).Value2 = 500
$ws.Cells.Item(329, 12).Value2 = 500
$ws.Cells.Item(329, 13).Value2 = 500
$ws.Cells.Item(329, 14).Value  = "$/paquete 5 unidades"
$ws.Cells.Item(329, 15).Value  = "Región Metropolitana"
$ws.Cells.Item(329, 16).Value2 = 100
$ws.Cells.Item(329, 17).Value2 = 5
$ws.Cells.Item(329, 18).Value  = "Hortaliza"

# Keep the date column formatted the same way as the rest of column D.
$ws.Range("D328:D329").NumberFormat = $ws.Range("D330").NumberFormat
